$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $Val)
    $r = $Sheet.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.ClearFormats()
}

Set-TextValue $ws "D2" "51.588.04"
$ws.Range("E2").Value = "  +1.01%  "
Set-TextValue $ws "D3" "3.024.73"
$ws.Range("E3").Value = "  +2.20%  "
Set-TextValue $ws "D4" "1.00"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws "D5" "379.75"
$ws.Range("E5").Value = "  -0.15%  "
Set-TextValue $ws "D6" "102.21"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.01%  "
Set-TextValue $ws "D10" "36.64"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  +0.86%  "
Set-TextValue $ws "D13" "3.503.33"
$ws.Range("E13").Value = "  +2.20%  "
Set-TextValue $ws "D14" "18.51"
$ws.Range("E14").Value = "  +0.55%  "
Set-TextValue $ws "D15" "7.73"
$ws.Range("E15").Value = "  -0.23%  "
Set-TextValue $ws "D16" "3.019.60"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("E17").Value = "  -3.63%  "
Set-TextValue $ws "D18" "10.62"
$ws.Range("E18").Value = "  -14.11%  "
Set-TextValue $ws "D19" "51.612.88"
$ws.Range("E19").Value = "  +0.95%  "
Set-TextValue $ws "D20" "3.10"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  +0.30%  "
Set-TextValue $ws "D22" "0.0₃0960"
$ws.Range("E22").Value = "  +0.20%  "
Set-TextValue $ws "D23" "69.92"
$ws.Range("E23").Value = "  +0.30%  "
Set-TextValue $ws "D24" "267.19"
$ws.Range("E24").Value = "  -0.65%  "
Set-TextValue $ws "D25" "3.15"
$ws.Range("E25").Value = "  -5.86%  "
Set-TextValue $ws "D26" "8.46"
$ws.Range("E26").Value = "  +6.36%  "
Set-TextValue $ws "D27" "7.50"
$ws.Range("E27").Value = "  +7.21%  "
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  -0.03%  "
Set-TextValue $ws "D30" "26.15"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -0.62%  "
Set-TextValue $ws "D32" "10.27"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D33" "34.00"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D34" "50.58"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D35" "0.0449"
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D36" "2.02"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("E37").Value = "  -0.16%  "
Set-TextValue $ws "D38" "3.31"
$ws.Range("E38").Value = "  +1.39%  "
Set-TextValue $ws "D39" "0.301"
$ws.Range("E39").Value = "  +14.79%  "
Set-TextValue $ws "D40" "17.06"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D41" "129.03"
$ws.Range("E41").Value = "  +4.13%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D42" "1.85"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  +5.15%  "
Set-TextValue $ws "D46" "21.56"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D47" "2.07"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D48" "2.43"
$ws.Range("E48").Value = "  +3.49%  "
Set-TextValue $ws "D49" "2.019.06"
$ws.Range("E49").Value = "  -3.45%  "
Set-TextValue $ws "D50" "3.324.37"
$ws.Range("E50").Value = "  +2.13%  "
Set-TextValue $ws "D51" "0.515"
$ws.Range("E51").Value = "  +5.54%  "
